$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.165.97'
$ws.Range("E2").Value = '  -1.45%  '
$ws.Range("D3").Value = '2.184.21'
$ws.Range("E3").Value = '  -2.44%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.624'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '67.80'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.36%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.574'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.88'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.36%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0939'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.98%  '
$ws.Range("B12").Value = 'Avalanche'
$ws.Range("C12").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '35.74'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -12.68%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.104'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.01%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.62%  '
$ws.Range("D15").Value = '2.510.91'
$ws.Range("E15").Value = '  -2.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.871'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.32'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.67%  '
$ws.Range("D18").Value = '2.156.66'
$ws.Range("E18").Value = '  -3.54%  '
$ws.Range("D19").Value = '41.084.06'
$ws.Range("E19").Value = '  -1.55%  '
$ws.Range("D20").Value = '0.0₃0953'
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.17'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.52'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +14.59%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("B27").Value = 'WEMIXToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.45'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.73%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.15'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.99%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.80'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.120'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.65'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0756'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.123'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.73%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.32'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.60'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.69'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0308'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.21'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.05'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.57'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '62.88'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.196'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.22'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.65'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.100'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.91%  '
$ws.Range("E49").Value = '  +0.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.56%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '91.88'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.40%  '
